# This script re-applies a set of full-row content swaps/rotations on the
# "Artfynd" worksheet. A handful of rows that describe the same field
# observation were re-ordered (e.g. because of a re-sort upstream); every
# cell in a given row moves together as a block to another row in the same
# small group of rows. The mapping below says: "row X should end up
# containing exactly what row Y currently contains".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (source row's current/original full content
# becomes the new content of the destination row)
$mapping = @{
    9  = 10
    10 = 9
    15 = 17
    16 = 15
    17 = 16
    32 = 33
    33 = 32
    36 = 37
    37 = 36
    38 = 39
    39 = 38
    43 = 44
    44 = 43
    48 = 49
    49 = 48
    52 = 53
    53 = 52
    54 = 55
    55 = 54
    56 = 57
    57 = 58
    58 = 59
    59 = 56
}

$lastCol = "AY"

# 1) Snapshot every involved row's full contents (A..AY) before mutating
#    anything, so that rotations (3- and 4-row cycles) don't clobber data
#    that a later row still needs to read.
$snapshot = @{}
foreach ($r in $mapping.Keys) {
    $addr = "A{0}:{1}{0}" -f $r, $lastCol
    $snapshot[$r] = $ws.Range($addr).Value2
}

# Also snapshot the date/time columns as their displayed Text, so we can
# restore them verbatim afterwards (Excel auto-converts strings such as
# "2026-01-25" or "09:07" into real date/time serials when they are written
# back through .Value2, which would change the cell's stored type).
$dateCols = @("Y", "Z", "AA", "AB")
$textSnapshot = @{}
foreach ($r in $mapping.Keys) {
    $textSnapshot[$r] = @{}
    foreach ($col in $dateCols) {
        $textSnapshot[$r][$col] = $ws.Range("$col$r").Text
    }
}

# 2) Write each destination row's new content from the snapshot of its
#    source row.
foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    $addr = "A{0}:{1}{0}" -f $r, $lastCol
    $ws.Range($addr).Value2 = $snapshot[$src]
}

# 3) Restore the date/time-looking text cells explicitly as plain text so
#    they keep the same literal representation they had in the source row
#    (rather than becoming Excel date/time serial numbers).
foreach ($r in $mapping.Keys) {
    $src = $mapping[$r]
    foreach ($col in $dateCols) {
        $cell = $ws.Range("$col$r")
        $text = $textSnapshot[$src][$col]
        if ([string]::IsNullOrEmpty($text)) {
            $cell.Value2 = ""
        } else {
            $cell.NumberFormat = "@"
            $cell.Value2 = $text
        }
    }
}
